$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1442.4
$ws.Range("I43").Value = 1000
$ws.Range("J43").Value = 1553
$ws.Range("K43").Value = 1000
$ws.Range("L43").Value = 1553
$ws.Range("M43").Value = -931
$ws.Range("N43").Value = -1691

$ws.Range("H64").Value = 3302.366
$ws.Range("I64").Value = 3071.2144
$ws.Range("J64").Value = 3422.2222
$ws.Range("K64").Value = 3071.2144
$ws.Range("L64").Value = 3422.2222
$ws.Range("M64").Value = -2823.2144
$ws.Range("N64").Value = -3918.2222

$ws.Range("H67").Value = 3302.366
$ws.Range("I67").Value = 3071.2144
$ws.Range("J67").Value = 3422.2222
$ws.Range("K67").Value = 3071.2144
$ws.Range("L67").Value = 3422.2222
$ws.Range("M67").Value = -2213.2144
$ws.Range("N67").Value = -5138.2222

$ws.Range("H76").Value = 4636.636
$ws.Range("I76").Value = 5300.5
$ws.Range("J76").Value = 3840
$ws.Range("K76").Value = 5300.5
$ws.Range("L76").Value = 3840
$ws.Range("M76").Value = -4985.5
$ws.Range("N76").Value = -4470

$ws.Range("H79").Value = 4636.636
$ws.Range("I79").Value = 5300.5
$ws.Range("J79").Value = 3840
$ws.Range("K79").Value = 5300.5
$ws.Range("L79").Value = 3840
$ws.Range("M79").Value = -4208.5
$ws.Range("N79").Value = -6024

$ws.Range("H86").Value = 2158.2307
$ws.Range("I86").Value = 2205.3
$ws.Range("J86").Value = 2001.3334
$ws.Range("K86").Value = 2205.3
$ws.Range("L86").Value = 2001.3334
$ws.Range("M86").Value = -1082.3
$ws.Range("N86").Value = -4247.3334

$ws.Range("H89").Value = 2158.2307
$ws.Range("I89").Value = 2205.3
$ws.Range("J89").Value = 2001.3334
$ws.Range("K89").Value = 11026.5
$ws.Range("L89").Value = 10006.667
$ws.Range("M89").Value = -5410.5
$ws.Range("N89").Value = -21238.667

$ws.Range("H121").Value = 1521.9231
$ws.Range("I121").Value = 1056.6666
$ws.Range("J121").Value = 1661.5
$ws.Range("K121").Value = 3169.9998
$ws.Range("L121").Value = 4984.5
$ws.Range("M121").Value = -1422.9998
$ws.Range("N121").Value = -8478.5

$ws.Range("H137").Value = 1761.6177
$ws.Range("I137").Value = 1682.3158
$ws.Range("J137").Value = 1862.0667
$ws.Range("K137").Value = 5046.9474
$ws.Range("L137").Value = 5586.2001
$ws.Range("M137").Value = -2496.9474
$ws.Range("N137").Value = -10686.2001

$ws.Range("H138").Value = 3237.9443
$ws.Range("I138").Value = 2329.543
$ws.Range("J138").Value = 4097.243
$ws.Range("K138").Value = 6988.629000000001
$ws.Range("L138").Value = 12291.729
$ws.Range("M138").Value = -1848.629000000001
$ws.Range("N138").Value = -22571.729

$ws.Range("H141").Value = 4055.681
$ws.Range("I141").Value = 1870.0465
$ws.Range("J141").Value = 27551.25
$ws.Range("K141").Value = 5610.139499999999
$ws.Range("L141").Value = 82653.75
$ws.Range("M141").Value = -430.1394999999993
$ws.Range("N141").Value = -93013.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1581.3617
$ws.Range("I132").Value = 1189.04
$ws.Range("J132").Value = 2027.1818
$ws.Range("K132").Value = 3567.12
$ws.Range("L132").Value = 6081.5454
$ws.Range("M132").Value = -1037.12
$ws.Range("N132").Value = -11141.5454

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1735.279
$ws.Range("I31").Value = 1308.1316
$ws.Range("J31").Value = 4981.6
$ws.Range("K31").Value = 1308.1316
$ws.Range("L31").Value = 4981.6
$ws.Range("M31").Value = -1013.1316
$ws.Range("N31").Value = -5571.6

$ws.Range("H34").Value = 1735.279
$ws.Range("I34").Value = 1308.1316
$ws.Range("J34").Value = 4981.6
$ws.Range("K34").Value = 1308.1316
$ws.Range("L34").Value = 4981.6
$ws.Range("M34").Value = -1106.1316
$ws.Range("N34").Value = -5385.6

$ws.Range("H134").Value = 1372.7273
$ws.Range("I134").Value = 1120.925
$ws.Range("J134").Value = 2044.2
$ws.Range("K134").Value = 3362.775
$ws.Range("L134").Value = 6132.6
$ws.Range("M134").Value = -827.7749999999996
$ws.Range("N134").Value = -11202.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 189.07692
$ws.Range("I7").Value = 119.77778
$ws.Range("J7").Value = 345
$ws.Range("K7").Value = 359.33334
$ws.Range("L7").Value = 1035
$ws.Range("M7").Value = -247.33334
$ws.Range("N7").Value = -1259

$ws.Range("H33").Value = 1238.3636
$ws.Range("I33").Value = 476.85715
$ws.Range("J33").Value = 2571
$ws.Range("K33").Value = 2861.1429
$ws.Range("L33").Value = 15426
$ws.Range("M33").Value = -2578.1429
$ws.Range("N33").Value = -15992

$ws.Range("H69").Value = 578
$ws.Range("I69").Value = 337.33334
$ws.Range("J69").Value = 1300
$ws.Range("K69").Value = 1012.00002
$ws.Range("L69").Value = 3900
$ws.Range("M69").Value = -201.0000200000001
$ws.Range("N69").Value = -5522

$ws.Range("H72").Value = 578
$ws.Range("I72").Value = 337.33334
$ws.Range("J72").Value = 1300
$ws.Range("K72").Value = 3036.00006
$ws.Range("L72").Value = 11700
$ws.Range("M72").Value = 1019.99994
$ws.Range("N72").Value = -19812

$ws.Range("H117").Value = 797.9091
$ws.Range("I117").Value = 561.3333
$ws.Range("J117").Value = 1081.8
$ws.Range("K117").Value = 1683.9999
$ws.Range("L117").Value = 3245.4
$ws.Range("M117").Value = 1758.0001
$ws.Range("N117").Value = -10129.4

$ws.Range("H121").Value = 1093
$ws.Range("I121").Value = 565
$ws.Range("J121").Value = 1726.6
$ws.Range("K121").Value = 1695
$ws.Range("L121").Value = 5179.799999999999
$ws.Range("M121").Value = -385
$ws.Range("N121").Value = -7799.799999999999

$ws.Range("H122").Value = 680.8
$ws.Range("I122").Value = 498.84616
$ws.Range("J122").Value = 877.9167
$ws.Range("K122").Value = 4489.61544
$ws.Range("L122").Value = 7901.2503
$ws.Range("M122").Value = -2039.61544
$ws.Range("N122").Value = -12801.2503

$ws.Range("H132").Value = 2350
$ws.Range("I132").Value = 1700
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 15300
$ws.Range("L132").Value = 27000
$ws.Range("M132").Value = -12770
$ws.Range("N132").Value = -32060

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2466.3635
$ws.Range("I122").Value = 1923.3334
$ws.Range("J122").Value = 4910
$ws.Range("K122").Value = 5770.0002
$ws.Range("L122").Value = 14730
$ws.Range("M122").Value = -3320.0002
$ws.Range("N122").Value = -19630

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4095.423
$ws.Range("I7").Value = 3426
$ws.Range("J7").Value = 4586.3335
$ws.Range("K7").Value = 3426
$ws.Range("L7").Value = 4586.3335
$ws.Range("M7").Value = -3314
$ws.Range("N7").Value = -4810.3335

$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("M87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("M90").ClearContents()

$ws.Range("H126").Value = 4095.423
$ws.Range("I126").Value = 3426
$ws.Range("J126").Value = 4586.3335
$ws.Range("K126").Value = 10278
$ws.Range("L126").Value = 13759.0005
$ws.Range("M126").Value = -7808
$ws.Range("N126").Value = -18699.0005

$ws.Range("H132").Value = 3396.6155
$ws.Range("I132").Value = 3054.6667
$ws.Range("J132").Value = 7500
$ws.Range("K132").Value = 9164.000100000001
$ws.Range("L132").Value = 22500
$ws.Range("M132").Value = -6634.000100000001
$ws.Range("N132").Value = -27560

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 22217.7
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 22217.7
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 22217.7
$ws.Range("N123").Value = -32017.7
